# Update pipeline results in the "Validacion_Rangos" worksheet.
# Refreshes Min_Observado / Max_Observado values (columns D and E) and the
# derived Estado_Rango / Analisis_Severidad text (columns F and G) for rows 2-8,
# reflecting the updated 'Presion_vapor_entrada' minimum range and the
# refreshed real-data path pipeline run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Validacion_Rangos")

# Row 2: 1GEV007CE
$ws.Range("D2").Value = 263.56
$ws.Range("E2").Value = 272.11
$ws.Range("F2").Value = "OK: Dentro de rango"
$ws.Range("G2").Value = "Cálculo pendiente (Métrica de severidad)"

# Row 3: 1FSRFTB504
$ws.Range("D3").Value = 774.38
$ws.Range("E3").Value = 788.41
$ws.Range("F3").Value = "OK: Dentro de rango"
$ws.Range("G3").Value = "Cálculo pendiente (Métrica de severidad)"

# Row 4: 1FSRTE502C
$ws.Range("D4").Value = 401.19
$ws.Range("E4").Value = 403.35
$ws.Range("F4").Value = "OK: Dentro de rango"
$ws.Range("G4").Value = "Cálculo pendiente (Métrica de severidad)"

# Row 5: 1FRSTE503A
$ws.Range("D5").Value = 535.14
$ws.Range("E5").Value = 538.88
$ws.Range("F5").Value = "OK: Dentro de rango"
$ws.Range("G5").Value = "Cálculo pendiente (Métrica de severidad)"

# Row 6: 1FSRPT501
$ws.Range("D6").Value = 171.45
$ws.Range("E6").Value = 174.16
$ws.Range("F6").Value = "OK: Dentro de rango"
$ws.Range("G6").Value = "Cálculo pendiente (Métrica de severidad)"

# Row 7: 1FSRPT504 (min observado now below expected minimum)
$ws.Range("D7").Value = 164.26
$ws.Range("E7").Value = 166.8
$ws.Range("F7").Value = "BAJO: Bajo el minimo esperado"
$ws.Range("G7").Value = "Cálculo pendiente (Métrica de severidad)"

# Row 8: 1FRSPT526
$ws.Range("D8").Value = 29.51
$ws.Range("E8").Value = 30.38
$ws.Range("F8").Value = "OK: Dentro de rango"
$ws.Range("G8").Value = "Cálculo pendiente (Métrica de severidad)"
